$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts Barcode/BC Added/etc. one column right)
$ws.Columns("F").Insert()

# New column header
$ws.Range("F1").Value = "Plate_Barcode"

# Excel auto-fit the new column to the header text when it was added
$ws.Columns("F").AutoFit()

# Leave the cursor where the author left it
$ws.Range("F4").Select() | Out-Null
